$wb = $excel.ActiveWorkbook
$ch11 = $wb.Worksheets.Item("Chapter 11")
$review = $wb.Worksheets.Item("Review")

# Chapter 11 review: the grader copied the "correct answer" column (B1:B15)
# and pasted it twice into the Review sheet as new columns AE and AF.

# 1) Paste the values into AE1:AE15
$ch11.Range("B1:B15").Copy()
$review.Range("AE1").PasteSpecial()

# 2) Paste the values again into AF1:AF15
$ch11.Range("B1:B15").Copy()
$review.Range("AF1").PasteSpecial()

# 3) The source rows 9 & 10 carry special (red) formatting in Chapter 11,
#    which becomes the plain/no-highlight style (style index 2, same as the
#    rest of the "right-hand" paired columns on the Review sheet) once
#    pasted. Re-apply that formatting (sourced from an existing cell that
#    already carries it) to AF9 / AF10 only - AE keeps the default style.
$review.Range("T1").Copy()
$review.Range("AF9").PasteSpecial(-4122)
$review.Range("AF10").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# 4) Restore the selections left behind by the edit.
$ch11.Range("B1:B15").Select()
$review.Range("V24").Select()
